$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.162.14"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.846.54"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "279.73"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5110"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3506"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.02"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06836"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.97"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8075"
$ws.Range("E12").Value = "  -5.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07781"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "1.849.01"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.099"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.58"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9990"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.20"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008077"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "26.199.67"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.775"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.210"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.385"
$ws.Range("E25").Value = "  +10.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.41"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.659"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.24"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.18"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.375"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.312"
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08746"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04912"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.172"
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7381"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.238"
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01860"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5178"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9648"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.25"
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.264"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.015"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4535"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1360"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.362"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.37"
$ws.Range("E49").Value = "  +1.72%  "

# Row 50/51 swap: Cronos <-> NEARProtocol with updated values
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.502"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05923"
$ws.Range("E51").Value = "  +0.32%  "
